$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at 569 (pushes the former rows 569:645 down to 570:646,
# carrying the style of row 569's cells, e.g. the date style on column D).
$ws.Rows("569:569").Insert()

# Populate the newly inserted row 569 with the new record's data.
$ws.Range("A569").Value = 8
$ws.Range("B569").Value = "Terminal La Palmera de La Serena"
$ws.Range("C569").Value = "Coquimbo"
$ws.Range("D569").Value = 45077
$ws.Range("E569").Value = 4
$ws.Range("F569").Value = 100114001
$ws.Range("G569").Value = "Papa"
$ws.Range("H569").Value = "Asterix"
$ws.Range("I569").Value = "1a nueva(o)"
$ws.Range("J569").Value = 2000
$ws.Range("K569").Value = 11800
$ws.Range("L569").Value = 12000
$ws.Range("M569").Value = 11900
$ws.Range("N569").Value = "$/saco 25 kilos"
$ws.Range("O569").Value = "Provincia de Melipilla"
$ws.Range("P569").Value = 476
$ws.Range("Q569").Value = 25
$ws.Range("R569").Value = "Hortaliza"
